$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition list) — update "想去人数" (interested count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2095
$ws1.Range("F4").Value = 866
$ws1.Range("F5").Value = 1275

# Sheet "全部类型" (all types combined list) — same three events repeated
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2095
$ws4.Range("F6").Value = 866
$ws4.Range("F7").Value = 1275
